# Insert a new row at position 213, shifting existing rows 213:336 down to 214:337,
# then populate the new row 213 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("213:213").Insert()

$ws.Cells.Item(213, 1).Value = 9
$ws.Cells.Item(213, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(213, 3).Value = "Metropolitana"
$ws.Cells.Item(213, 4).Value = 44873
$ws.Cells.Item(213, 4).NumberFormat = $ws.Cells.Item(214, 4).NumberFormat
$ws.Cells.Item(213, 5).Value = 13
$ws.Cells.Item(213, 6).Value = 300000001
$ws.Cells.Item(213, 7).Value = "Rabanito"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 14000
$ws.Cells.Item(213, 11).Value = 3000
$ws.Cells.Item(213, 12).Value = 4000
$ws.Cells.Item(213, 13).Value = 3643
$ws.Cells.Item(213, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(213, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(213, 16).Value = 36
$ws.Cells.Item(213, 17).Value = 100
$ws.Cells.Item(213, 18).Value = "Hortaliza"
